$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their original Text storage type (avoid Excel
# auto-converting numeric-looking strings like "176.40" into numbers).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.951.06'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +4.37%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.248.38'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +2.51%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.72'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +3.56%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '176.40'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +2.88%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.41%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.246.39'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +2.47%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +4.93%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.72'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.64%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.408'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +3.26%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.812.42'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +2.59%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.84'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.21%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.935.10'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +4.36%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000168'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +3.73%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.246.75'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +2.58%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.82'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.84%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.31'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +2.36%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '369.23'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +4.94%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.48'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +4.37%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.07%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '70.45'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +2.07%  '
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = 'Polygon'
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.507'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.35%  '
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = 'WrappedeETH'
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.380.01'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +2.32%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.55%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.79'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +3.98%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.98%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.18%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +5.45%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.63'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.75%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '22.49'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +2.08%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.19%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '172.81'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +9.85%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.23'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +2.99%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.77'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.67%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.51'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +4.87%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.852'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +6.93%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +10.21%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '26.83'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +2.81%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.58'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +2.44%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.40'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +6.16%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.714.38'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +2.52%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.30'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +3.81%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +4.60%  '
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'Bittensor'
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '339.22'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +4.44%  '
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = 'Hedera'
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0673'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +3.27%  '
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '24.51'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +4.09%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0279'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +3.52%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +2.55%  '
